$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each cell as literal text (NumberFormat '@' prevents Excel from
# reinterpreting numeric-looking / date-looking strings such as '229.00',
# '1.00', or '37.840.24' as numbers).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.840.24'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.085.13'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.74%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.73'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.09%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.89'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.98%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.99%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.83%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.392.78'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.61%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.76'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.22'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.770'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.42%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.085.56'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.55%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.740.53'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.99%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.19'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.36'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0835'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.00'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.39'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.59%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '170.61'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.140'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +7.17%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.92%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.42'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.42%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.54'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.24%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.67%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.69'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.73'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.30%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.37%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.52'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.45'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.06%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.93%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.14%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.95%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0994'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.14%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '98.71'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.39'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +4.75%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.465.96'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.73%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.63%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.07'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.52%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.99'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +5.46%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.40'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.37%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.59%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.277.33'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.50%  '
